$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - 1.1 User Stories
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 9

# Row 6 - 1.3 ER Diagram
$ws.Range("B6").Value = 2
$ws.Range("H6").Value = 9

# Row 7 - 1.3 Cost Estimation
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("F7").Value = 3
$ws.Range("H7").Value = 11

# Row 8 - 1.4 Proposed Screens
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 3
$ws.Range("H8").Value = 12

# Row 9 - 1.5 Schedule Management
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 11

# Row 10 - Sprint 2 - User Development (total)
$ws.Range("H10").Value = 52

# Row 11 - 2.1 User UI development
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 2
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 16

# Row 12 - 2.2 User Database Design
$ws.Range("B12").Value = 4
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 15

# Row 13 - 2.3 User Implementation
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = 15

# Row 14 - 2.4 User Deploy
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 16

# Row 15 - Sprint 3 - Admin Development (total)
$ws.Range("G15").ClearContents()
$ws.Range("H15").Value = 62

# Row 16 - 3.1 Admin UI development
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 1

# Row 17 - 3.2 Admin Database Design
$ws.Range("C17").Value = 1
$ws.Range("G17").Value = 2

# Row 18 - 3.3 Admin Implementation
$ws.Range("D18").Value = 4
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 4

# Row 19 - 3.4 Admin Deploy
$ws.Range("D19").Value = 2
$ws.Range("G19").Value = 2

# Update selection to I16
$ws.Range("I16").Select()
